# Resubmission of project 1 - corrected errors in analysis and report.
# The "Data" sheet gains a third measurement column (Differences = Incongruent
# - Congruent) and the summary statistics block (previously in columns D:K)
# is rebuilt two columns further right (F:M) to make room for it. The
# "mean/SEM" statistic is recomputed directly from the new differences
# column instead of being derived algebraically from the two sample SEMs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------
# 1. Clear out the old D1:K25 summary block so we can rebuild it in the
#    new F:M location without leftovers.
# ---------------------------------------------------------------------
$ws.Range("D1:K25").ClearContents()

# ---------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Congruent"
$ws.Range("B1").Value = "Incongruent"
$ws.Range("C1").Value = "Differences"
$ws.Range("D1").Value = "Bins"
$ws.Range("F1").Value = "Congruent"
$ws.Range("I1").Value = "Incongruent"

# ---------------------------------------------------------------------
# 3. New "Differences" column C = Incongruent - Congruent.
#    C2 is entered on its own; C3:C25 is filled as one relative-formula
#    range (mirrors how the original workbook's columns were authored -
#    one seed cell plus a fill-down range).
# ---------------------------------------------------------------------
$ws.Range("C2").Formula = "=B2-A2"
$ws.Range("C3:C25").Formula = "=B3-A3"
$ws.Range("C2:C27").NumberFormat = "General"

# ---------------------------------------------------------------------
# 4. Histogram bin edges, moved from C2:C5 to D2:D5.
# ---------------------------------------------------------------------
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 20
$ws.Range("D4").Value = 30
$ws.Range("D5").Value = 40

# ---------------------------------------------------------------------
# 5. Congruent descriptive stats: labels in F, values in G (was D/E).
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "x"
$ws.Range("G2").Formula = "=AVERAGE(A2:A25)"
$ws.Range("F3").Value = "s"
$ws.Range("G3").Formula = "=_xlfn.STDEV.S(A2:A25)"
$ws.Range("F4").Value = "n"
$ws.Range("G4").Formula = "=COUNTA(A2:A25)"

# ---------------------------------------------------------------------
# 6. Incongruent descriptive stats: labels in I, values in J (was G/H).
# ---------------------------------------------------------------------
$ws.Range("I2").Value = "x"
$ws.Range("J2").Formula = "=AVERAGE(B2:B25)"
$ws.Range("I3").Value = "s"
$ws.Range("J3").Formula = "=_xlfn.STDEV.S(B2:B25)"
$ws.Range("I4").Value = "n"
$ws.Range("J4").Formula = "=COUNTA(B2:B25)"

# ---------------------------------------------------------------------
# 7. Paired-difference test block: labels in L, values in M (was J/K).
#    d    = mean of the differences column (was H2-E2)
#    s(d) = sample stdev of the differences column (was se(d)=SQRT(s1^2+s2^2))
#    se(d)= s(d)/sqrt(n)            (re-labelled "se(d)", was "sem")
#    df, t-statistic, t-critical, ci low/high keep the same shape, just
#    re-pointed at the new M-column cells.
# ---------------------------------------------------------------------
$ws.Range("L2").Value = "d"
$ws.Range("M2").Formula = "=AVERAGE(C2:C25)"

$ws.Range("L3").Value = "s(d)"
$ws.Range("M3").Formula = "=_xlfn.STDEV.S(C2:C25)"

$ws.Range("L4").Value = "se(d)"
$ws.Range("M4").Formula = "=M3/SQRT(G4)"

$ws.Range("L5").Value = "df"
$ws.Range("M5").Formula = "=G4-1"

$ws.Range("L6").Value = "t-statistic"
$ws.Range("M6").Formula = "=M2/(M3/SQRT(G4))"

$ws.Range("L7").Value = "t-critical for a = 0.05"
$ws.Range("M7").Value = 2.069

$ws.Range("L8").Value = "ci low"
$ws.Range("M8").Formula = "=`$M`$2-`$M`$7*`$M`$4"

$ws.Range("L9").Value = "ci high"
$ws.Range("M9").Formula = "=`$M`$2+`$M`$7*`$M`$4"

# ---------------------------------------------------------------------
# 8. A couple of trailing blank-but-formatted cells under the
#    Differences column (rows 26-27), matching the extended used range.
# ---------------------------------------------------------------------
$ws.Range("C26").NumberFormat = "General"
$ws.Range("C27").NumberFormat = "General"

# ---------------------------------------------------------------------
# 9. Column width: column A keeps its bestFit width, and the wide
#    "ci high" label column moves from J to M.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.625
$ws.Columns.Item(11).ColumnWidth = 17.5

# ---------------------------------------------------------------------
# 10. View state: zoomed to 85%, current selection on H14.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("H14").Select()
